{"js": "// The document ends with a paragraph that holds the \"_GoBack\" bookmark and\n// (before this edit) the diary entry text \"\u591a\u4e91\uff0c\u4eca\u5929\u662f\u56fd\u5e86\u8282\". This edit:\n//   1. Moves that entry's text out into its own new paragraph.\n//   2. Adds a new date/weather diary-entry pair right after it:\n//        \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\"\n//        \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\"\n//   3. Leaves the original trailing paragraph empty, keeping only the\n//      bookmark in place.\nconst body = context.document.body;\n\n// Locate the final paragraph (the one carrying the _GoBack bookmark).\nconst bmRange = context.document.getBookmarkRange(\"_GoBack\");\nconst tailPara = bmRange.paragraphs.getFirst();\ntailPara.load(\"text\");\n\nconst paras = body.paragraphs;\nparas.load(\"items\");\n\nawait context.sync();\n\nconst tailText = tailPara.text;\n\n// Find the paragraph immediately before the tail paragraph; the new entries\n// get inserted \"after\" it so each one inherits its (eastAsia) run-hint\n// formatting, matching how Word splits off a new paragraph when typing\n// right before an existing one.\nconst items = paras.items;\nconst tailIndex = items.length - 1;\nlet anchor = items[tailIndex - 1];\n\nconst newEntries = [\n  tailText,\n  \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\",\n  \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\",\n];\n\nfor (const entryText of newEntries) {\n  anchor = anchor.insertParagraph(entryText, Word.InsertLocation.after);\n}\nawait context.sync();\n\n// Strip the old text out of the (now shifted) tail paragraph, leaving the\n// bookmark untouched.\nconst results = body.search(tailText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nconst target = results.items[results.items.length - 1];\ntarget.delete();\nawait context.sync();\n", "ps1": "# The document ends with a paragraph that holds the \"_GoBack\" bookmark and\n# (before this edit) the diary entry text \"\u591a\u4e91\uff0c\u4eca\u5929\u662f\u56fd\u5e86\u8282\". This edit:\n#   1. Moves that entry's text out into its own new paragraph.\n#   2. Adds a new date/weather diary-entry pair right after it:\n#        \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\"\n#        \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\"\n#   3. Leaves the original trailing paragraph empty, keeping only the\n#      bookmark in place.\n$d = $word.ActiveDocument\n\n# Locate the final paragraph (the one carrying the _GoBack bookmark).\n$bm = $d.Bookmarks(\"_GoBack\")\n$tailIndex = $bm.Range.Paragraphs.First.Index\n$tailText = $d.Paragraphs.Item($tailIndex).Range.Text -replace \"[\\r\\a\\f]+$\", \"\"\n\n# New diary entries to insert right before the (soon to be emptied) tail\n# paragraph: first the text that used to live in the tail paragraph, then a\n# new date line, then a new weather/description line.\n$newEntries = @(\n    $tailText,\n    \"2022\u5e746\u67087\u65e5\u661f\u671f\u4e8c\",\n    \"\u6674\uff0c\u4eca\u5929\u662f\u9ad8\u8003\u7b2c\u4e00\u5929\uff0c\u4e0a\u5348\u8003\u8bed\u6587\uff0c\u4e0b\u5348\u8003\u6570\u5b66\u3002\"\n)\n\n$anchorIndex = $tailIndex - 1\nforeach ($entryText in $newEntries) {\n    $anchor = $d.Paragraphs.Item($anchorIndex)\n    $anchor.Range.InsertParagraphAfter()\n    $anchorIndex = $anchorIndex + 1\n    $d.Paragraphs.Item($anchorIndex).Range.Text = $entryText\n}\n\n# Strip the old text out of the (now shifted) tail paragraph, leaving the\n# bookmark untouched.\n$newTailIndex = $d.Paragraphs.Count\n$tailPara = $d.Paragraphs.Item($newTailIndex)\n$tailRange = $d.Range($tailPara.Range.Start, $tailPara.Range.End)\n$tailRange.Find.Execute($tailText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n"}
